$wb = $excel.ActiveWorkbook

function Update-Cell {
    param($ws, $cell, $oldValue, $newValue)
    $current = $ws.Range($cell).Value
    $ws.Range($cell).Value = $newValue
}

$ws1 = $wb.Worksheets.Item("展览")
Update-Cell $ws1 "F3" 86 87
Update-Cell $ws1 "F4" 9796 9805
Update-Cell $ws1 "F9" 431 432
Update-Cell $ws1 "F11" 213 215
Update-Cell $ws1 "F12" 475 476
Update-Cell $ws1 "F13" 12432 12435
Update-Cell $ws1 "F21" 184 185
Update-Cell $ws1 "F24" 2739 2740
Update-Cell $ws1 "F26" 87 88
Update-Cell $ws1 "F30" 1053 1055
Update-Cell $ws1 "F31" 4226 4228
Update-Cell $ws1 "F32" 3724 3726
Update-Cell $ws1 "F33" 728 730
Update-Cell $ws1 "F34" 2637 2638
Update-Cell $ws1 "F35" 3068 3070
Update-Cell $ws1 "F38" 200 201
Update-Cell $ws1 "F42" 459 460
Update-Cell $ws1 "F43" 595 598

$ws2 = $wb.Worksheets.Item("演出")
Update-Cell $ws2 "F13" 46 47
Update-Cell $ws2 "F20" 80 81

$ws4 = $wb.Worksheets.Item("全部类型")
Update-Cell $ws4 "F4" 86 87
Update-Cell $ws4 "F5" 9796 9805
Update-Cell $ws4 "F11" 431 432
Update-Cell $ws4 "F13" 213 215
Update-Cell $ws4 "F14" 475 476
Update-Cell $ws4 "F15" 12432 12435
Update-Cell $ws4 "F21" 184 185
Update-Cell $ws4 "F24" 2739 2740
Update-Cell $ws4 "F26" 87 88
Update-Cell $ws4 "F29" 1053 1055
Update-Cell $ws4 "F30" 4226 4228
Update-Cell $ws4 "F31" 3724 3726
Update-Cell $ws4 "F32" 728 730
Update-Cell $ws4 "F33" 2637 2638
Update-Cell $ws4 "F34" 3068 3070
Update-Cell $ws4 "F37" 200 201
Update-Cell $ws4 "F41" 459 460
Update-Cell $ws4 "F43" 595 598
